$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 289, pushing the old row 289
# (Paine / 1a (guarda)) down to row 291, unchanged.
$ws.Rows.Item(289).Insert()
$ws.Rows.Item(289).Insert()

# Row 284: Camote, 1a (guarda), now priced 310/800/830/815, origin Provincia de Talca
$ws.Cells.Item(284, 4).Value = 44448
$ws.Cells.Item(284, 8).Value = "Camote"
$ws.Cells.Item(284, 9).Value = "1a (guarda)"
$ws.Cells.Item(284, 10).Value = 310
$ws.Cells.Item(284, 11).Value = 800
$ws.Cells.Item(284, 12).Value = 830
$ws.Cells.Item(284, 13).Value = 815
$ws.Cells.Item(284, 15).Value = "Provincia de Talca"
$ws.Cells.Item(284, 16).Value = 815

# Row 285: Paine, 1a (guarda), priced 80/550/550/550
$ws.Cells.Item(285, 4).Value = 44448
$ws.Cells.Item(285, 8).Value = "Paine"
$ws.Cells.Item(285, 9).Value = "1a (guarda)"
$ws.Cells.Item(285, 10).Value = 80
$ws.Cells.Item(285, 11).Value = 550
$ws.Cells.Item(285, 12).Value = 550
$ws.Cells.Item(285, 13).Value = 550
$ws.Cells.Item(285, 16).Value = 550

# Row 286: Camote, 1a nueva(o), priced 80/1300/1300/1300, origin Peru
$ws.Cells.Item(286, 4).Value = 44167
$ws.Cells.Item(286, 8).Value = "Camote"
$ws.Cells.Item(286, 9).Value = "1a nueva(o)"
$ws.Cells.Item(286, 10).Value = 80
$ws.Cells.Item(286, 11).Value = 1300
$ws.Cells.Item(286, 12).Value = 1300
$ws.Cells.Item(286, 13).Value = 1300
$ws.Cells.Item(286, 15).Value = "Perú"
$ws.Cells.Item(286, 16).Value = 1300

# Row 287: Camote, 1a nueva(o), priced 170/400/450/426
$ws.Cells.Item(287, 4).Value = 44238
$ws.Cells.Item(287, 8).Value = "Camote"
$ws.Cells.Item(287, 9).Value = "1a nueva(o)"
$ws.Cells.Item(287, 10).Value = 170
$ws.Cells.Item(287, 11).Value = 400
$ws.Cells.Item(287, 12).Value = 450
$ws.Cells.Item(287, 13).Value = 426
$ws.Cells.Item(287, 16).Value = 426

# Row 288: Paine, 1a nueva(o), priced 80/350/350/350
$ws.Cells.Item(288, 4).Value = 44238
$ws.Cells.Item(288, 8).Value = "Paine"
$ws.Cells.Item(288, 9).Value = "1a nueva(o)"
$ws.Cells.Item(288, 10).Value = 80
$ws.Cells.Item(288, 11).Value = 350
$ws.Cells.Item(288, 12).Value = 350
$ws.Cells.Item(288, 13).Value = 350
$ws.Cells.Item(288, 16).Value = 350

# Row 289 (new): Camote, 1a (guarda), priced 280/400/420/413
$ws.Cells.Item(289, 1).Value = 3
$ws.Cells.Item(289, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(289, 3).Value = "Coquimbo"
$ws.Cells.Item(289, 4).Value = 44399
$ws.Cells.Item(289, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(289, 5).Value = 5
$ws.Cells.Item(289, 6).Value = 100112045
$ws.Cells.Item(289, 7).Value = "Zapallo"
$ws.Cells.Item(289, 8).Value = "Camote"
$ws.Cells.Item(289, 9).Value = "1a (guarda)"
$ws.Cells.Item(289, 10).Value = 280
$ws.Cells.Item(289, 11).Value = 400
$ws.Cells.Item(289, 12).Value = 420
$ws.Cells.Item(289, 13).Value = 413
$ws.Cells.Item(289, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(289, 15).Value = "Provincia de Talca"
$ws.Cells.Item(289, 16).Value = 413
$ws.Cells.Item(289, 17).Value = 1
$ws.Cells.Item(289, 18).Value = "Hortaliza"

# Row 290 (new): Camote, 1a (guarda), priced 185/400/420/410
$ws.Cells.Item(290, 1).Value = 3
$ws.Cells.Item(290, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(290, 3).Value = "Coquimbo"
$ws.Cells.Item(290, 4).Value = 44400
$ws.Cells.Item(290, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(290, 5).Value = 5
$ws.Cells.Item(290, 6).Value = 100112045
$ws.Cells.Item(290, 7).Value = "Zapallo"
$ws.Cells.Item(290, 8).Value = "Camote"
$ws.Cells.Item(290, 9).Value = "1a (guarda)"
$ws.Cells.Item(290, 10).Value = 185
$ws.Cells.Item(290, 11).Value = 400
$ws.Cells.Item(290, 12).Value = 420
$ws.Cells.Item(290, 13).Value = 410
$ws.Cells.Item(290, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(290, 15).Value = "Provincia de Talca"
$ws.Cells.Item(290, 16).Value = 410
$ws.Cells.Item(290, 17).Value = 1
$ws.Cells.Item(290, 18).Value = "Hortaliza"

# Row 291 already holds the former row-289 data (Paine / 1a (guarda) / 330)
# after the row inserts above, so no further edits are required there.

$wb.Save()
